$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new field "igm_iga_agree" right after "igm_igg_agree" (old row 29), ---
# --- before "antibody_agree" (old row 30) ---
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 2).Value = "igm_iga_agree"
$ws.Cells.Item(30, 3).Value = "string"
$ws.Cells.Item(30, 4).Value = "Agreement between igm_iga_result and antibody_truth."
$ws.Rows.Item(30).RowHeight = 17

# --- Insert new field "igm_iga_result" right after "igm_igg_result" (row 17), ---
# --- before "control" (old row 18) ---
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 2).Value = "igm_iga_result"
$ws.Cells.Item(18, 3).Value = "string"
$ws.Cells.Item(18, 4).Value = "The test result for qualitative detection of (IgM / IgA) combined antibodies."
$ws.Rows.Item(18).RowHeight = 17

# --- Update the view selection to D19 ---
$ws.Range("D19").Select()
